$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - Dolar (USD)
$ws.Range("E2").Value = "5,54"
$ws.Range("F2").Value = "11 de jun., 18:46 UTC ·"

# Row 3 - Euro (EUR)
$ws.Range("E3").Value = "6,37"
$ws.Range("F3").Value = "11 de jun., 18:46 UTC ·"

# Row 4 - Real (BRL)
$ws.Range("F4").Value = "11 de jun., 18:46 UTC ·"

# Row 5 - Libra Esterlina (GBP)
$ws.Range("E5").Value = "7,51"
$ws.Range("F5").Value = "11 de jun., 18:46 UTC ·"

# Row 6 - Iene (JPY)
$ws.Range("F6").Value = "11 de jun., 18:47 UTC ·"

# Row 7 - Franco Suico (CHF)
$ws.Range("F7").Value = "11 de jun., 18:46 UTC ·"

# Row 8 - Dolar Australiano (AUD)
$ws.Range("F8").Value = "11 de jun., 18:46 UTC ·"

# Row 9 - Peso Mexicano (MXN)
$ws.Range("F9").Value = "11 de jun., 18:46 UTC ·"

# Row 10 - Dolar Canadiano (CAD)
$ws.Range("E10").Value = "4,06"
$ws.Range("F10").Value = "11 de jun., 18:46 UTC ·"

# Row 11 - Dolar de Hong Kong (HKD)
$ws.Range("F11").Value = "11 de jun., 18:46 UTC ·"

# Row 12 - Yuan Chines (CNY)
$ws.Range("F12").Value = "11 de jun., 18:45 UTC ·"

# Row 13 - Rupia Indiana (INR)
$ws.Range("F13").Value = "11 de jun., 18:45 UTC ·"

# Row 14 - Peso Chileno (CLP) - no change

# Row 15 - Peso Argentino (ARS)
$ws.Range("F15").Value = "11 de jun., 18:46 UTC ·"

# Row 16 - Peso Colombiano (COP)
$ws.Range("F16").Value = "11 de jun., 18:46 UTC ·"

# Row 17 - Rupia Russa (RUB)
$ws.Range("F17").Value = "11 de jun., 18:46 UTC ·"

# Row 18 - Riyal Saudi (SAR)
$ws.Range("F18").Value = "11 de jun., 18:46 UTC ·"

# Row 19 - Dolar de Singapura (SGD)
$ws.Range("F19").Value = "11 de jun., 18:46 UTC ·"

# Row 20 - Peso Filipino (PHP)
$ws.Range("F20").Value = "11 de jun., 18:46 UTC ·"

# Row 21 - Yuan de Taiwan (TWD)
$ws.Range("E21").Value = "4,16"
$ws.Range("F21").Value = "11 de jun., 18:47 UTC ·"

# Row 22 - Dinar Iraquiano (IQD)
$ws.Range("F22").Value = "11 de jun., 18:46 UTC ·"

# Row 23 - Rupia Sri Lanka (LKR)
$ws.Range("F23").Value = "11 de jun., 18:45 UTC ·"

# Row 24 - Yuan Chines (CNY)
$ws.Range("F24").Value = "11 de jun., 18:45 UTC ·"

# Row 25 - Won Sul-Coreano (KRW)
$ws.Range("F25").Value = "11 de jun., 18:47 UTC ·"
